$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '57.626.30'
$ws.Range("E2").Value = '  -1.93%  '

$ws.Range("D3").Value = '3.099.85'
$ws.Range("E3").Value = '  -2.41%  '

$ws.Range("E4").Value = '  -0.03%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '528.77'
$ws.Range("D5").Style = $ws.Range("B5").Style
$ws.Range("E5").Value = '  -1.27%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '137.75'
$ws.Range("D6").Style = $ws.Range("B6").Style
$ws.Range("E6").Value = '  -3.52%  '

$ws.Range("E7").Value = '  -0.18%  '

$ws.Range("D8").Value = '3.097.94'
$ws.Range("E8").Value = '  -2.49%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.466'
$ws.Range("D9").Style = $ws.Range("B9").Style

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '7.27'
$ws.Range("D10").Style = $ws.Range("B10").Style
$ws.Range("E10").Value = '  +0.29%  '

$ws.Range("E11").Value = '  -3.28%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.407'
$ws.Range("D12").Style = $ws.Range("B12").Style
$ws.Range("E12").Value = '  +1.67%  '

$ws.Range("E13").Value = '  +1.78%  '

$ws.Range("D14").Value = '3.631.81'
$ws.Range("E14").Value = '  -2.53%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '25.42'
$ws.Range("D15").Style = $ws.Range("B15").Style
$ws.Range("E15").Value = '  -2.11%  '

$ws.Range("E16").Value = '  -3.37%  '

$ws.Range("D17").Value = '57.644.95'
$ws.Range("E17").Value = '  -2.07%  '

$ws.Range("D18").Value = '3.088.24'
$ws.Range("E18").Value = '  -3.00%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '5.93'
$ws.Range("D19").Style = $ws.Range("B19").Style
$ws.Range("E19").Value = '  -4.24%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '12.55'
$ws.Range("D20").Style = $ws.Range("B20").Style
$ws.Range("E20").Value = '  -3.13%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '7.92'
$ws.Range("D21").Style = $ws.Range("B21").Style
$ws.Range("E21").Value = '  -2.15%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '350.24'
$ws.Range("D22").Style = $ws.Range("B22").Style
$ws.Range("E22").Value = '  -2.15%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '1.00'
$ws.Range("D23").Style = $ws.Range("B23").Style
$ws.Range("E23").Value = '  -0.03%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '68.79'
$ws.Range("D24").Style = $ws.Range("B24").Style
$ws.Range("E24").Value = '  +0.55%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.502'
$ws.Range("D25").Style = $ws.Range("B25").Style
$ws.Range("E25").Value = '  -2.53%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.167'
$ws.Range("D26").Style = $ws.Range("B26").Style
$ws.Range("E26").Value = '  -2.28%  '

$ws.Range("E27").Value = '  +0.17%  '

$ws.Range("D28").Value = '0.0₃0867'
$ws.Range("E28").Value = '  -9.87%  '

$ws.Range("E29").Value = '  -0.05%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '7.14'
$ws.Range("D30").Style = $ws.Range("B30").Style
$ws.Range("E30").Value = '  -5.43%  '

$ws.Range("E31").Value = '  -2.53%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '6.00'
$ws.Range("D32").Style = $ws.Range("B32").Style
$ws.Range("E32").Value = '  -8.63%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '21.12'
$ws.Range("D33").Style = $ws.Range("B33").Style
$ws.Range("E33").Value = '  -1.34%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.89'
$ws.Range("D34").Style = $ws.Range("B34").Style
$ws.Range("E34").Value = '  -0.65%  '

$ws.Range("B35").Value = 'Fetch.AI'
$ws.Range("C35").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.14'
$ws.Range("D35").Style = $ws.Range("B35").Style
$ws.Range("E35").Value = '  -6.92%  '

$ws.Range("B36").Value = 'Monero'
$ws.Range("C36").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '158.81'
$ws.Range("D36").Style = $ws.Range("B36").Style
$ws.Range("E36").Value = '  +0.59%  '

$ws.Range("E37").Value = '  -3.50%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '25.68'
$ws.Range("D38").Style = $ws.Range("B38").Style
$ws.Range("E38").Value = '  -3.23%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.25'
$ws.Range("D39").Style = $ws.Range("B39").Style
$ws.Range("E39").Value = '  -5.19%  '

$ws.Range("B40").Value = 'Stacks'
$ws.Range("C40").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.64'
$ws.Range("D40").Style = $ws.Range("B40").Style
$ws.Range("E40").Value = '  -1.47%  '

$ws.Range("B41").Value = 'Hedera'
$ws.Range("C41").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0666'
$ws.Range("D41").Style = $ws.Range("B41").Style
$ws.Range("E41").Value = '  -1.92%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '4.03'
$ws.Range("D42").Style = $ws.Range("B42").Style
$ws.Range("E42").Value = '  -0.53%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.693'
$ws.Range("D43").Style = $ws.Range("B43").Style
$ws.Range("E43").Value = '  -2.43%  '

$ws.Range("D44").Value = '2.396.31'
$ws.Range("E44").Value = '  +2.14%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '36.96'
$ws.Range("D45").Style = $ws.Range("B45").Style
$ws.Range("E45").Value = '  +0.02%  '

$ws.Range("E46").Value = '  -0.13%  '

$ws.Range("D47").Value = '3.138.82'
$ws.Range("E47").Value = '  -2.44%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0264'
$ws.Range("D48").Style = $ws.Range("B48").Style
$ws.Range("E48").Value = '  -4.00%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '6.03'
$ws.Range("D49").Style = $ws.Range("B49").Style
$ws.Range("E49").Value = '  -0.91%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.953'
$ws.Range("D50").Style = $ws.Range("B50").Style
$ws.Range("E50").Value = '  -6.91%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '19.59'
$ws.Range("D51").Style = $ws.Range("B51").Style
$ws.Range("E51").Value = '  -5.51%  '
